$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray leading-space run (rFonts Helvetica / color 464048)
#    that sits just before "What you should have:" (first FOR-block's
#    heading). There is a second, visually-identical run before "What will
#    make you thrive:" later in the doc that must be left untouched, so we
#    locate the specific paragraph by its trailing text instead of relying
#    on a global Find/Replace.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*What you should have:*") {
        $pStart = $p.Range.Start
        $firstChar = $d.Range($pStart, $pStart + 1)
        if ($firstChar.Font.Color -eq 4735046) {
            $firstChar.Delete()
        }
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Insert a new "+++IF specificCompetence != null+++" paragraph right
#    after the first "+++END-FOR comp+++" paragraph (core competences loop)
#    and right before the "What will make you thrive:" heading.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $n; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($t -eq "+++END-FOR comp+++") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ne -1) {
    $anchor = $d.Paragraphs.Item($targetIndex)
    $anchor.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($targetIndex + 1)

    $ifXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="172B4D"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="172B4D"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">+++IF </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="172B4D"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>specificCompetence</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="172B4D"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> !</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="172B4D"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>= null+++</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $newPara.Range.InsertXML($ifXml)
}

# ---------------------------------------------------------------------------
# 3) Add <w:rStyle w:val="normaltextrun"/> as the first child of the
#    paragraph-mark rPr for the second "+++END-FOR comp+++" paragraph
#    (specific competences loop, trailing spaces run in Calibri), then
#    append a new "+++END-IF+++" paragraph right after it.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $n; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "+++END-FOR comp+++*") {
        $targetIndex = $i
    }
}

if ($targetIndex -ne -1) {
    $anchor = $d.Paragraphs.Item($targetIndex)

    # Clone the paragraph's own OOXML (preserving every existing attribute)
    # and splice in the missing <w:rStyle> on the paragraph-mark rPr, which
    # is the very first <w:rPr> emitted (the pPr's one).
    $full = $anchor.Range.WordOpenXML
    if ($full -match '(?s)<w:body>(.*)</w:body>') {
        $bodyFrag = $matches[1]
        $rPrIdx = $bodyFrag.IndexOf("<w:rPr>")
        if ($rPrIdx -ge 0) {
            $bodyFrag = $bodyFrag.Substring(0, $rPrIdx + 7) + '<w:rStyle w:val="normaltextrun"/>' + $bodyFrag.Substring($rPrIdx + 7)
            $wrapped = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $bodyFrag + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
            $anchor.Range.InsertXML($wrapped)
        }
    }

    $anchor = $d.Paragraphs.Item($targetIndex)
    $anchor.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($targetIndex + 1)

    $endIfXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="normaltextrun"/><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>+++END-IF+++</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $newPara.Range.InsertXML($endIfXml)
}
